# Highlight quantitative impact metrics (percentages, dollar amounts, etc.)
# throughout the resume body with hybrid bold + color (#2C3E50) formatting.

function Set-MetricHighlight {
    param($doc, [string]$paraMatch, [string[]]$metrics)

    $paras = $doc.Paragraphs
    $target = $null
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like $paraMatch) {
            $target = $p
            break
        }
    }
    if ($null -eq $target) {
        Write-Output "PARAGRAPH NOT FOUND: $paraMatch"
        return
    }

    foreach ($metric in $metrics) {
        $r = $target.Range
        $r.Find.ClearFormatting()
        $found = $r.Find.Execute($metric, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $r.Font.Bold = $true
            $r.Font.Color = 5258796
        } else {
            Write-Output "METRIC NOT FOUND: '$metric' in paragraph matching '$paraMatch'"
        }
    }
}

$d = $word.ActiveDocument

# Partner - Siege Analytics bullets
Set-MetricHighlight $d "*Discovered systematic race coding errors*" @("23%", "64%")
Set-MetricHighlight $d "*Utilized advanced sampling methods*" @("±4.2%", "±2.1%", "71%", "87%")
Set-MetricHighlight $d "*Trigonometric algorithm for boundary estimation*" @("73.5%", "`$4.7M")
Set-MetricHighlight $d "*Built real-time FEC analysis systems*" @("`$2")

# Data Products Manager - Helm/Murmuration bullet
Set-MetricHighlight $d "*Modernized legacy ETL processes*" @("57%")

# KEY ACHIEVEMENTS AND IMPACT bullets
Set-MetricHighlight $d "*Algorithmic innovation: Pioneered trigonometric*" @("73.5%")
Set-MetricHighlight $d "*`$4.7M savings enabled nonprofit access*" @("`$4.7M")
Set-MetricHighlight $d "*178% accuracy improvement*" @("178%")

Write-Output "Metric highlighting complete"
